# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row on both the zh-cn and de-de sheets, simulating
# a freshly regenerated handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 00:51:09"
$wsZhCn.Range("H2").Value = "2016-03-19 00:51:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 00:51:12"
$wsDeDe.Range("H2").Value = "2016-03-19 00:51:37"
